$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.235.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.852.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4602"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3705"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07293"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8864"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.02"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07818"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.863.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.384"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.522"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.003"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008919"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.003"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.249.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.108"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.079.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.87%  "
$ws.Range("E25").Value = "  +5.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.055"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.059"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08822"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.096"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7676"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.167"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.497"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.746"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +11.66%  "
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05253"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.948"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.053"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5115"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1631"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.382"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4793"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.76%  "
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.641"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "101.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06214"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.39%  "
